$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.169.54'
$ws.Range("E2").Value = '  -0.51%  '

$ws.Range("D3").Value = '2.502.34'
$ws.Range("E3").Value = '  +1.59%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '540.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.73%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.76'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.91%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.16%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.572'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.67%  '

$ws.Range("D9").Value = '2.526.78'
$ws.Range("E9").Value = '  +2.19%  '

$ws.Range("E10").Value = '  +1.02%  '

$ws.Range("E11").Value = '  +0.34%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.59'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.98%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.354'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.67%  '

$ws.Range("D14").Value = '2.945.01'
$ws.Range("E14").Value = '  +1.24%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.53'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.16%  '

$ws.Range("D16").Value = '59.050.86'
$ws.Range("E16").Value = '  -0.74%  '

$ws.Range("E17").Value = '  +1.22%  '

$ws.Range("D18").Value = '2.521.49'
$ws.Range("E18").Value = '  +1.18%  '

$ws.Range("E19").Value = '  +0.60%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.28'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.64%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '325.03'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.23%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.86%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.78'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.71%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '62.19'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.40%  '

$ws.Range("E25").Value = '  -4.57%  '

$ws.Range("E26").Value = '  +0.89%  '

$ws.Range("D27").Value = '2.620.88'
$ws.Range("E27").Value = '  +2.13%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.994'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.99%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.79'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.81%  '

$ws.Range("D30").Value = '0.0₃0775'
$ws.Range("E30").Value = '  +0.24%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.82'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.20%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.67'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.77%  '

$ws.Range("E33").Value = '  -5.32%  '

$ws.Range("E34").Value = '  -0.26%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.44'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.54%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '156.46'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.89%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.66'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.38%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.34'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.63%  '

$ws.Range("E39").Value = '  -9.47%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.71'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.79%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.91'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.54%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '295.90'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.01%  '

$ws.Range("E43").Value = '  -0.59%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.820'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.58%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.996'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.30%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.600'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.64%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.77'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.42%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0930'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.04%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '122.60'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.84%  '

$ws.Range("E50").Value = '  +0.04%  '

$ws.Range("E51").Value = '  -0.28%  '
